$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    # Swap column D and E (4 and 5)
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal

    # Swap column F and H (6 and 8)
    $fVal = $ws.Cells.Item($r, 6).Value2
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 6).Value2 = $hVal
    $ws.Cells.Item($r, 8).Value2 = $fVal
}
